$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 25350.889
$ws.Range("J17").Value = 25350.889
$ws.Range("L17").Value = 76052.667
$ws.Range("N17").Value = -76388.667
$ws.Range("H70").Value = 1168.5714
$ws.Range("I70").Value = 955.8
$ws.Range("J70").Value = 1286.7778
$ws.Range("K70").Value = 2867.4
$ws.Range("L70").Value = 3860.3334
$ws.Range("M70").Value = -2597.4
$ws.Range("N70").Value = -4400.3334
$ws.Range("H73").Value = 1168.5714
$ws.Range("I73").Value = 955.8
$ws.Range("J73").Value = 1286.7778
$ws.Range("K73").Value = 2867.4
$ws.Range("L73").Value = 3860.3334
$ws.Range("M73").Value = -1931.4
$ws.Range("N73").Value = -5732.3334

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3343.8125
$ws.Range("I2").Value = 2885.077
$ws.Range("K2").Value = 2885.077
$ws.Range("M2").Value = -2772.077
$ws.Range("H63").Value = 2648.5881
$ws.Range("J63").Value = 3332.6667
$ws.Range("L63").Value = 3332.6667
$ws.Range("N63").Value = -4704.6667
$ws.Range("H66").Value = 2648.5881
$ws.Range("J66").Value = 3332.6667
$ws.Range("L66").Value = 16663.3335
$ws.Range("N66").Value = -23527.3335
$ws.Range("H80").Value = 19996.666
$ws.Range("J80").Value = 19996.666
$ws.Range("L80").Value = 19996.666
$ws.Range("N80").Value = -21992.666
$ws.Range("H83").Value = 19996.666
$ws.Range("J83").Value = 19996.666
$ws.Range("L83").Value = 59989.99800000001
$ws.Range("N83").Value = -69973.99800000001
$ws.Range("H116").Value = 3343.8125
$ws.Range("I116").Value = 2885.077
$ws.Range("K116").Value = 2885.077
$ws.Range("M116").Value = -591.0770000000002
$ws.Range("H132").Value = 5444.3584
$ws.Range("I132").Value = 6270.075
$ws.Range("K132").Value = 18810.225
$ws.Range("M132").Value = -16280.225

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3343.8125
$ws.Range("I3").Value = 2885.077
$ws.Range("K3").Value = 2885.077
$ws.Range("M3").Value = -2771.077
$ws.Range("H134").Value = 3371.75
$ws.Range("I134").Value = 3344
$ws.Range("K134").Value = 10032
$ws.Range("M134").Value = -7497

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2449.162
$ws.Range("I31").Value = 1879.1111
$ws.Range("J31").Value = 2989.2104
$ws.Range("K31").Value = 1879.1111
$ws.Range("L31").Value = 2989.2104
$ws.Range("M31").Value = -1584.1111
$ws.Range("N31").Value = -3579.2104
$ws.Range("H34").Value = 2449.162
$ws.Range("I34").Value = 1879.1111
$ws.Range("J34").Value = 2989.2104
$ws.Range("K34").Value = 1879.1111
$ws.Range("L34").Value = 2989.2104
$ws.Range("M34").Value = -1677.1111
$ws.Range("N34").Value = -3393.2104
$ws.Range("H82").Value = 25181
$ws.Range("J82").Value = 25181
$ws.Range("L82").Value = 25181
$ws.Range("N82").Value = -25903
$ws.Range("H85").Value = 25181
$ws.Range("J85").Value = 25181
$ws.Range("L85").Value = 25181
$ws.Range("N85").Value = -27677
$ws.Range("H107").Value = 1476.421
$ws.Range("I107").Value = 1011.375
$ws.Range("J107").Value = 1814.6364
$ws.Range("K107").Value = 1011.375
$ws.Range("L107").Value = 1814.6364
$ws.Range("M107").Value = 908.625
$ws.Range("N107").Value = -5654.6364
$ws.Range("H122").Value = 2920.4092
$ws.Range("I122").Value = 2847
$ws.Range("J122").Value = 3077.7144
$ws.Range("K122").Value = 8541
$ws.Range("L122").Value = 9233.143199999999
$ws.Range("M122").Value = -6091
$ws.Range("N122").Value = -14133.1432
$ws.Range("H132").Value = 2057.923
$ws.Range("I132").Value = 2057.923
$ws.Range("K132").Value = 6173.768999999999
$ws.Range("M132").Value = -3643.768999999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1560.0588
$ws.Range("I5").Value = 1451.2858
$ws.Range("J5").Value = 1636.2
$ws.Range("K5").Value = 4353.857400000001
$ws.Range("L5").Value = 4908.6
$ws.Range("M5").Value = -4241.857400000001
$ws.Range("N5").Value = -5132.6
$ws.Range("H18").Value = 812.4
$ws.Range("I18").Value = 812.4
$ws.Range("K18").Value = 2437.2
$ws.Range("M18").Value = -2268.2
$ws.Range("H33").Value = 256.15384
$ws.Range("I33").Value = 92.15385000000001
$ws.Range("J33").Value = 420.15384
$ws.Range("K33").Value = 552.9231
$ws.Range("L33").Value = 2520.92304
$ws.Range("M33").Value = -269.9231
$ws.Range("N33").Value = -3086.92304
$ws.Range("H38").Value = 110
$ws.Range("I38").Value = 95
$ws.Range("J38").Value = 125
$ws.Range("K38").Value = 285
$ws.Range("L38").Value = 375
$ws.Range("M38").Value = 62
$ws.Range("N38").Value = -1069
$ws.Range("H116").Value = 2978.5386
$ws.Range("I116").Value = 1524.7778
$ws.Range("J116").Value = 6249.5
$ws.Range("K116").Value = 4574.3334
$ws.Range("L116").Value = 18748.5
$ws.Range("M116").Value = -1132.3334
$ws.Range("N116").Value = -25632.5
$ws.Range("H130").Value = 1999.3334
$ws.Range("I130").Value = 1999.3334
$ws.Range("K130").Value = 5998.0002
$ws.Range("M130").Value = -978.0002000000004
$ws.Range("H135").Value = 1560.0588
$ws.Range("I135").Value = 1451.2858
$ws.Range("J135").Value = 1636.2
$ws.Range("K135").Value = 13061.5722
$ws.Range("L135").Value = 14725.8
$ws.Range("M135").Value = -10526.5722
$ws.Range("N135").Value = -19795.8

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8188.636
$ws.Range("I80").Value = 11523.583
$ws.Range("J80").Value = 4186.7
$ws.Range("K80").Value = 11523.583
$ws.Range("L80").Value = 4186.7
$ws.Range("M80").Value = -10525.583
$ws.Range("N80").Value = -6182.7
$ws.Range("H83").Value = 8188.636
$ws.Range("I83").Value = 11523.583
$ws.Range("J83").Value = 4186.7
$ws.Range("K83").Value = 57617.915
$ws.Range("L83").Value = 20933.5
$ws.Range("M83").Value = -52625.915
$ws.Range("N83").Value = -30917.5
$ws.Range("H123").Value = 59568.5
$ws.Range("J123").Value = 59568.5
$ws.Range("L123").Value = 59568.5
$ws.Range("N123").Value = -64468.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 96900
$ws.Range("J125").Value = 96900
$ws.Range("L125").Value = 96900
$ws.Range("N125").Value = -101820
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 30000
$ws.Range("I129").Value = 30000
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 30000
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -25000
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 97000
$ws.Range("J131").Value = 97000
$ws.Range("L131").Value = 97000
$ws.Range("N131").Value = -107080

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 356.6316
$ws.Range("I55").Value = 303.54544
$ws.Range("K55").Value = 303.54544
$ws.Range("M55").Value = -130.54544
$ws.Range("H93").Value = 27541.924
$ws.Range("I93").Value = 2095.4285
$ws.Range("K93").Value = 2095.4285
$ws.Range("M93").Value = -847.4285

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1858
$ws.Range("I81").Value = 1972.625
$ws.Range("J81").Value = 941
$ws.Range("K81").Value = 3945.25
$ws.Range("L81").Value = 1882
$ws.Range("M81").Value = -2884.25
$ws.Range("N81").Value = -4004
$ws.Range("H84").Value = 1858
$ws.Range("I84").Value = 1972.625
$ws.Range("J84").Value = 941
$ws.Range("K84").Value = 19726.25
$ws.Range("L84").Value = 9410
$ws.Range("M84").Value = -14422.25
$ws.Range("N84").Value = -20018
